$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 2.425633666666667
$ws.Range("N2").Value = 7.276901000000001
$ws.Range("O2").Value = 0.0662600404061536
$ws.Range("P2").Value = 0.06626004040615362
$ws.Range("Q2").Value = 0.5242449274756666
$ws.Range("R2").Value = 4.718204347281
$ws.Range("S2").Value = 0.03402239974736252
$ws.Range("T2").Value = 0.03402239974736253

# Row 3
$ws.Range("O3").Value = 0.4234968256437875
$ws.Range("P3").Value = 0.4234968256437876
$ws.Range("S3").Value = 0.2174519998097362
$ws.Range("T3").Value = 0.2174519998097363

# Row 4
$ws.Range("M4").Value = 18.67887366666667
$ws.Range("N4").Value = 56.03662100000001
$ws.Range("O4").Value = 0.5102431339500588
$ws.Range("P4").Value = 0.5102431339500588
$ws.Range("Q4").Value = 4.037008928955667
$ws.Range("R4").Value = 36.33308036060101
$ws.Range("S4").Value = 0.2619934392612254
$ws.Range("T4").Value = 0.2619934392612254

# Row 5
$ws.Range("M5").Value = 2.425633666666667
$ws.Range("N5").Value = 7.276901000000001
$ws.Range("O5").Value = 0.0662600404061536
$ws.Range("P5").Value = 0.06626004040615362
$ws.Range("Q5").Value = 0.4967439015075556
$ws.Range("R5").Value = 4.470695113568
$ws.Range("S5").Value = 0.03223764065879108
$ws.Range("T5").Value = 0.03223764065879109

# Row 6
$ws.Range("O6").Value = 0.4234968256437875
$ws.Range("P6").Value = 0.4234968256437876
$ws.Range("S6").Value = 0.2060448258340513
$ws.Range("T6").Value = 0.2060448258340514

# Row 7
$ws.Range("M7").Value = 18.67887366666667
$ws.Range("N7").Value = 56.03662100000001
$ws.Range("O7").Value = 0.5102431339500588
$ws.Range("P7").Value = 0.5102431339500588
$ws.Range("Q7").Value = 3.825234085614224
$ws.Range("R7").Value = 34.42710677052801
$ws.Range("S7").Value = 0.2482496946888335
$ws.Range("T7").Value = 0.2482496946888335
